$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (currently sitting right after
#    the word "major").
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2) Tidy up the double space between "Gallifrey." and "The Doctor
#    travels" -> single space (part of the re-wording captured by the
#    diff).
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Gallifrey.  The Doctor travels", $true, $false, $false, $false, $false, $true, 1, $false, "Gallifrey. The Doctor travels", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Split the big paragraph into three paragraphs (with a blank
#    paragraph between each), in document order from the back so the
#    earlier Find targets remain valid/unique.
# ---------------------------------------------------------------------

# 3a) split right before "The Doctor is clever, mad, and passionate"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "The Doctor is clever, mad, and passionate"
$find.Execute() | Out-Null
$rng = $find.Parent.Duplicate
$rng.Collapse(1)
$rng.InsertParagraphBefore()

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "The Doctor is clever, mad, and passionate"
$find.Execute() | Out-Null
$rng = $find.Parent.Duplicate
$rng.Collapse(1)
$rng.InsertParagraphBefore()

# 3b) split right before "The Doctor travels with companions"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "The Doctor travels with companions"
$find.Execute() | Out-Null
$rng = $find.Parent.Duplicate
$rng.Collapse(1)
$rng.InsertParagraphBefore()

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "The Doctor travels with companions"
$find.Execute() | Out-Null
$rng = $find.Parent.Duplicate
$rng.Collapse(1)
$rng.InsertParagraphBefore()

# ---------------------------------------------------------------------
# 4) Split the run "...Time Lord from the Planet Gallifrey..." into
#    "...Time L" / "ord from the Planet Gallifrey. " / "The Doctor
#    travels through..." by toggling (and restoring) Bold across the
#    two split boundaries -- this forces the engine to break runs at
#    those exact points without altering the visible formatting.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "ord from the Planet Gallifrey. "
$find.Execute() | Out-Null
$splitRng = $find.Parent.Duplicate
$splitRng.Font.Bold = 1
$splitRng.Font.Bold = 0

# ---------------------------------------------------------------------
# 5) Re-create the "_GoBack" bookmark at the very start of the
#    (now-first) paragraph of this block, spanning through to the end
#    of the block's final paragraph.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Doctor Who, a British television show"
$find.Execute() | Out-Null
$startPos = $find.Parent.Start

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "Jodie Whittaker, the first female Doctor"
$find2.Execute() | Out-Null
$endPos = $find2.Parent.End

$bmRange = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "edit complete"
